$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New A/B values for rows 3-22 (row 2 stays unchanged)
$data = @(
    @(56, "W931101108060"),
    @(53, "J931101109013"),
    @(17, "M931252909052"),
    @(32, "X931412020027"),
    @(42, "M931321110016"),
    @(41, "T931252910046"),
    @(39, "Y931321110015"),
    @(3,  "Q931325208064"),
    @(24, "Y931321309029"),
    @(19, "E931100609021"),
    @(14, "U931412020025"),
    @(28, "B931383814058"),
    @(55, "D931100609028"),
    @(34, "R931100609011"),
    @(36, "Q931100609020"),
    @(15, "M931252110020"),
    @(46, "P931252710020"),
    @(47, "Q931101109046"),
    @(12, "R931325310022"),
    @(33, "V931252909047")
)

# Apply the same cell formatting (border/bold/alignment) used by the
# existing column-A data cells to the new rows (15-22) before filling
# in values.
$ws.Range("A2").Copy()
$ws.Range("A15:A22").PasteSpecial(-4122)

$row = 3
foreach ($pair in $data) {
    $ws.Cells.Item($row, 1).Value = $pair[0]
    $ws.Cells.Item($row, 2).Value = $pair[1]
    $row = $row + 1
}
